# Fruta / hortaliza, semanal
# Insert a new weekly record at row 434 (shifts the existing rows 434-476
# down to 435-477) in the "Hortaliza, Feria Lagunitas de Puerto Montt -
# Zapallo italiano" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 434..end down by one to make room for the new record.
$ws.Rows.Item(434).Insert()

# Populate the newly inserted row 434 with the new weekly data point.
$ws.Cells.Item(434, 1).Value = 4
$ws.Cells.Item(434, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(434, 3).Value = 'Los Lagos'
$ws.Cells.Item(434, 4).Value = 45212
$ws.Cells.Item(434, 5).Value = 10
$ws.Cells.Item(434, 6).Value = 100112032
$ws.Cells.Item(434, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(434, 8).Value = 'Sin especificar'
$ws.Cells.Item(434, 9).Value = 'Primera'
$ws.Cells.Item(434, 10).Value = 250
$ws.Cells.Item(434, 11).Value = 22000
$ws.Cells.Item(434, 12).Value = 22000
$ws.Cells.Item(434, 13).Value = 22000
$ws.Cells.Item(434, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(434, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(434, 16).Value = 440
$ws.Cells.Item(434, 17).Value = 50
$ws.Cells.Item(434, 18).Value = 'Hortaliza'
